# Weekly update: a new week's price record for Espinaca (Vega Modelo de Temuco)
# is inserted at the top of the data block (row 93), pushing the existing
# rows 93-121 down by one (to 94-122).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 93; this shifts rows 93:121 -> 94:122
# and keeps all existing data/formatting intact.
$ws.Rows("93:93").Insert()

# Populate the newly inserted row 93 with the new week's record.
$ws.Cells.Item(93, 1).Value = 10
$ws.Cells.Item(93, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(93, 3).Value = "La Araucanía"
$ws.Cells.Item(93, 4).Value = 44588
$ws.Cells.Item(93, 5).Value = 9
$ws.Cells.Item(93, 6).Value = 100112012
$ws.Cells.Item(93, 7).Value = "Espinaca"
$ws.Cells.Item(93, 8).Value = "Sin especificar"
$ws.Cells.Item(93, 9).Value = "Primera"
$ws.Cells.Item(93, 10).Value = 75
$ws.Cells.Item(93, 11).Value = 13000
$ws.Cells.Item(93, 12).Value = 14000
$ws.Cells.Item(93, 13).Value = 13533
$ws.Cells.Item(93, 14).Value = "`$/docena de atados"
$ws.Cells.Item(93, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(93, 16).Value = 4511
$ws.Cells.Item(93, 17).Value = 3
$ws.Cells.Item(93, 18).Value = "Hortaliza"
